# Aggiunta colonna "Id Materiale" mancante
# -------------------------------------------------------------------------
# Sheet "DB MATERIALI - Dati": insert a new column D ("ID MATERIALE"),
# shifting CATEGORIA..NOTE ULTERIORI one column to the right (D..K -> E..L).
# Sheet "DB MATERIALI - Conf": insert a new column E describing the new
# "ID MATERIALE" field, shifting CATEGORIA..NOTE ULTERIORI one column to
# the right (E..L -> F..M), and fill in the field-definition rows for it.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsDati = $wb.Worksheets.Item("DB MATERIALI - Dati")
$wsConf = $wb.Worksheets.Item("DB MATERIALI - Conf")

# --- Sheet "DB MATERIALI - Dati" -----------------------------------------
# Insert a blank column before D (CATEGORIA), pushing everything right.
$wsDati.Range("D1").EntireColumn.Insert()
$wsDati.Range("D1").Value = "ID MATERIALE"
$wsDati.Range("D1").ColumnWidth = 14.5

# --- Sheet "DB MATERIALI - Conf" -----------------------------------------
# Insert a blank column before E (CATEGORIA), pushing everything right.
$wsConf.Range("E1").EntireColumn.Insert()

$wsConf.Range("E1").Value = "ID MATERIALE"
$wsConf.Range("E2").Value = 3
$wsConf.Range("E3").Value = "ID MAT."
$wsConf.Range("E4").Value = "Testo"
$wsConf.Range("E6").Value = 255
$wsConf.Range("E7").Value = "No"
$wsConf.Range("E9").Value = "No"
$wsConf.Range("E10").Value = "Sì"
$wsConf.Range("E12").Value = "No"
